$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.414.97"
$ws.Range("E2").Value = "  -2.67%  "

$ws.Range("D3").Value = "2.213.69"
$ws.Range("E3").Value = "  -2.71%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "107.82"
$ws.Range("E5").Value = "  -12.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "294.79"
$ws.Range("E6").Value = "  +10.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -3.66%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("E9").Value = "  -4.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.63"
$ws.Range("E10").Value = "  -9.35%  "

$ws.Range("E11").Value = "  -4.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.53"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.76"
$ws.Range("E13").Value = "  -5.47%  "

$ws.Range("E14").Value = "  -3.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.958"
$ws.Range("E15").Value = "  +5.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.90"
$ws.Range("E16").Value = "  -3.47%  "

$ws.Range("D17").Value = "2.546.56"
$ws.Range("E17").Value = "  -2.72%  "

$ws.Range("D18").Value = "2.230.28"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").Value = "42.336.01"
$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  +5.99%  "

$ws.Range("E21").Value = "  -5.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.57"
$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.44"
$ws.Range("E23").Value = "  +18.82%  "

$ws.Range("E24").Value = "  -6.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "227.57"
$ws.Range("E25").Value = "  -3.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.00"
$ws.Range("E26").Value = "  -5.12%  "

$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.56"
$ws.Range("E28").Value = "  -3.59%  "

$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.47"
$ws.Range("E31").Value = "  -9.14%  "

$ws.Range("E32").Value = "  -5.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "173.46"
$ws.Range("E33").Value = "  +0.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.83"
$ws.Range("E34").Value = "  -4.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0884"
$ws.Range("E35").Value = "  -3.46%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.50"
$ws.Range("E36").Value = "  -4.01%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.99"
$ws.Range("E37").Value = "  +7.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("E38").Value = "  +1.35%  "

$ws.Range("E39").Value = "  -3.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0362"
$ws.Range("E40").Value = "  -3.78%  "

$ws.Range("E41").Value = "  -3.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.43"
$ws.Range("E42").Value = "  -4.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.85"
$ws.Range("E43").Value = "  -4.03%  "

$ws.Range("E44").Value = "  -4.11%  "

$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.64"
$ws.Range("E46").Value = "  -9.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.30"
$ws.Range("E47").Value = "  -5.96%  "

$ws.Range("E48").Value = "  -4.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.31"
$ws.Range("E49").Value = "  +3.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.05"
$ws.Range("E50").Value = "  +1.40%  "

$ws.Range("E51").Value = "  -1.45%  "
